{"js": "// RBA v2.5 - Atualizacao da Tela\n// Replace placeholder tokens (various casings of \"tre\"/\"tere\") with the\n// new placeholder tokens (various casings of \"qwer\"/\"qewr\"), both in the\n// main document body and in the page header, matching each occurrence in\n// document order.\n\n// --- Main document body: the single bold \"TERE\" inside \"A TERE, vem por...\" ---\nconst body = context.document.body;\nconst bodyResults = body.search(\"TERE\", { matchCase: true, matchWholeWord: false });\nbodyResults.load(\"items\");\nawait context.sync();\n\nif (bodyResults.items.length > 0) {\n  bodyResults.items[0].insertText(\"QWER\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// --- Header: replace each occurrence, in document order. ---\nconst sections = context.document.sections;\nsections.load(\"items\");\nawait context.sync();\n\nconst header = sections.items[0].getHeader(\"primary\");\n\nasync function replaceOrdered(term, replacements) {\n  const results = header.search(term, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length && i < replacements.length; i++) {\n    results.items[i].insertText(replacements[i], Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n// \"DIRETORIA DE ENSINO REGIAO TRE\" -> \"...QWER\"\nawait replaceOrdered(\"TRE\", [\"QWER\"]);\n// \"TERE - DEP.\" -> \"QWER - DEP.\"\nawait replaceOrdered(\"TERE\", [\"QWER\"]);\n// \"Tre, no Tre - Tre - Tre - Tre\" -> \"Qwer, no Qwer - Qewr - Qewr - Qwer\"\nawait replaceOrdered(\"Tre\", [\"Qwer\", \"Qwer\", \"Qewr\", \"Qewr\", \"Qwer\"]);\n// \"CEP: tre ... Tel: tre\" / \"Email: tre\" -> \"qwer\"\nawait replaceOrdered(\"tre\", [\"qwer\", \"qwer\", \"qwer\"]);\n", "ps1": "# RBA v2.5 - Atualizacao da Tela\n# Replace placeholder tokens (various casings of \"tre\"/\"tere\") with the\n# new placeholder tokens (various casings of \"qwer\"/\"qewr\"), both in the\n# main document body and in the page header, matching each occurrence in\n# document order.\n\n$d = $word.ActiveDocument\n\n# --- Main document body: the single bold \"TERE\" inside \"A TERE, vem por...\" ---\n$bodyRange = $d.Content\n$bodyRange.Find.Execute(\"TERE\", $true, $false, $false, $false, $false, $true, 1, $false, \"QWER\", 1) | Out-Null\n\n# --- Header: replace each occurrence, in document order, one at a time. ---\n# Re-fetching a fresh header Range for every call keeps Find positioned\n# correctly as earlier matches in the header are updated.\n$sec = $d.Sections.Item(1)\n\nfunction Replace-InHeader([string]$search, [string]$replace) {\n    $r = $sec.Headers.Item(1).Range\n    $r.Find.Execute($search, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 1) | Out-Null\n}\n\n# \"DIRETORIA DE ENSINO REGIAO TRE\" -> \"...QWER\"\nReplace-InHeader \"TRE\" \"QWER\"\n# \"TERE - DEP.\" -> \"QWER - DEP.\"\nReplace-InHeader \"TERE\" \"QWER\"\n# \"Tre, no Tre - Tre - Tre - Tre\" -> \"Qwer, no Qwer - Qewr - Qewr - Qwer\"\nReplace-InHeader \"Tre\" \"Qwer\"\nReplace-InHeader \"Tre\" \"Qwer\"\nReplace-InHeader \"Tre\" \"Qewr\"\nReplace-InHeader \"Tre\" \"Qewr\"\nReplace-InHeader \"Tre\" \"Qwer\"\n# \"CEP: tre ... Tel: tre\" / \"Email: tre\" -> \"qwer\"\nReplace-InHeader \"tre\" \"qwer\"\nReplace-InHeader \"tre\" \"qwer\"\nReplace-InHeader \"tre\" \"qwer\"\n"}
